$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.389.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.97%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.842.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  +1.28%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'315.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.13%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.012"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.19%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4746"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +1.75%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3705"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.58%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07472"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +1.31%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8863"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +1.93%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.52"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.70%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.857.93"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.16%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.07367"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +4.39%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.489"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.51%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +2.00%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'6.581"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.29%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'0.000008868"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +2.06%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.012"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.22%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.93%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'27.416.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.90%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.359"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.40%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +1.58%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.072.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.29%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'152.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.37%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.88%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.184"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.47%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.284"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.61%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'118.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.14%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08969"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.44%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.7627"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'1.181"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +1.54%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.46%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'2.948"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.62%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +1.28%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.107"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.07%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.05375"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.81%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01964"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.23%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.997"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.93%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.323"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.83%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.5363"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.01%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'2.384"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.63%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D45").Value = "'8.560"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.74%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4986"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.16%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +1.32%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.33%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'105.15"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.41%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'1.682"
$ws.Range("D50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.79%  "
$ws.Range("E51").Style = "Normal"
